$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.903.74"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.549.84"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0588"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "1.771.26"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "1.556.84"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "26.912.37"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("D34").Value = "1.411.10"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.965"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.529"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "1.685.28"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0517"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +5.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
